# Rewrite of example 11 unit tests
# Update the "Scenario" (column F) and "Given-When-Then (Description)" (column H)
# cells of the ATDD Scenarios table to reflect the renamed event subscribers /
# new trigger & check procedure names. Downstream formula columns
# (ATDD Format / Code Format / ATDD.TestScriptor Format) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATDD Scenarios")

# --- "Sales Document" posting scenarios ---------------------------------
$ws.Range("F3").Value  = "Check failure CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Rows.Item(3).RowHeight = 45.75
$ws.Range("H5").Value  = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Rows.Item(5).RowHeight = 30
$ws.Range("F7").Value  = "Check success CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Rows.Item(7).RowHeight = 45
$ws.Range("H9").Value  = "Trigger CheckLookupvalueExistsOnSalesHeader Sales Posting"
$ws.Rows.Item(9).RowHeight = 30

# --- "Warehouse Shipment" posting scenarios -----------------------------
$ws.Range("F11").Value = "Check failure CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Rows.Item(11).RowHeight = 45
$ws.Range("H13").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Rows.Item(13).RowHeight = 30
$ws.Range("F15").Value = "Check success CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Rows.Item(15).RowHeight = 45
$ws.Range("H17").Value = "Trigger CheckLookupvalueExistsOnSalesHeader Whse. Posting"
$ws.Rows.Item(17).RowHeight = 30

# --- Inheritance - Sales Document / Customer ----------------------------
$ws.Range("F20").Value = "Check InheritLookupValueFromCustomer"
$ws.Rows.Item(20).RowHeight = 30.75
$ws.Range("H23").Value = "Trigger InheritLookupValueFromCustomer"
$ws.Rows.Item(23).RowHeight = 16.5

# --- Inheritance - Contact ------------------------------------------------
$ws.Range("F25").Value = "Check ApplyLookupValueFromCustomerTemplate from Contact"
$ws.Rows.Item(25).RowHeight = 30
$ws.Range("H28").Value = "Trigger ApplyLookupValueFromCustomerTemplate"
$ws.Rows.Item(28).RowHeight = 30

# --- Inheritance - Customer Templates --------------------------------------
$ws.Range("F30").Value = "Check ApplyLookupValueFromCustomerTemplate"
$ws.Rows.Item(30).RowHeight = 30
$ws.Range("H33").Value = "Trigger ApplyLookupValueFromCustomerTemplate"
$ws.Rows.Item(33).RowHeight = 30

# --- Inheritance - Warehouse Shipment ---------------------------------------
$ws.Range("F36").Value = "Check InheritLookupValueFromSalesHeader"
$ws.Rows.Item(36).RowHeight = 30.75
$ws.Range("H39").Value = "Trigger InheritLookupValueFromSalesHeader"
$ws.Rows.Item(39).RowHeight = 30
